$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header title cells (A1/B1/C1) ---
# Net visible text changes from '16.5.1.1a.' to '16.5.1.1a ' (dot -> space) in all 3 languages.
$ws.Range("A1").Value = '16.5.1.1a "Аткаруу бийлигинин мамлекеттик органдарындагы жана жергиликтүү өз алдынча башкаруу органдарындагы коррупциянын деңгээли жөнүндө жеке түшүнүк" индекси'
$ws.Range("B1").Value = '16.5.1.1a Индекс "Личное представление об уровне коррупции в государственных органах исполнительной власти и органах местного самоуправления'''''
$ws.Range("C1").Value = '16.5.1.1a Index "Personal views about the level of corruption in executive government authorities and local government'''''

# --- Add new column I (year 2020 data), copying formatting from column H ---
$ws.Range("H4:H14").Copy()
$ws.Range("I4:I14").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("I4").Value = 2020

$ws.Range("I5:I14").NumberFormat = "0.0"
$ws.Range("I5").Value = 12.3
$ws.Range("I6").Value = 40.3
$ws.Range("I7").Value = 36.2
$ws.Range("I8").Value = 44.3
$ws.Range("I9").Value = 36
$ws.Range("I10").Value = 2.7
$ws.Range("I11").Value = 32.9
$ws.Range("I12").Value = 11.3
$ws.Range("I13").Value = -18.2
$ws.Range("I14").Value = 33

# --- Update selection to reflect the saved cursor position ---
$ws.Range("F16").Select()
